# Auto-generated edit script for 신세계푸드.xlsx (error solve ifrs list)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("P2").ClearContents()
$ws.Range("Q2").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
$ws.Range("AF2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AG2").Value = 750
$ws.Range("AH2").Value = 0.7
$ws.Range("AJ2").Value = 3872480

# Row 3
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("D3").Value = 9064
$ws.Range("E3").Value = 87
$ws.Range("F3").Value = 87
$ws.Range("G3").Value = 79
$ws.Range("H3").Value = 67
$ws.Range("I3").Value = 67
$ws.Range("K3").Value = 5033
$ws.Range("L3").Value = 2161
$ws.Range("M3").Value = 2871
$ws.Range("N3").Value = 2871
$ws.Range("P3").Value = 194
$ws.Range("Q3").Value = 127
$ws.Range("R3").Value = -1124
$ws.Range("S3").Value = 964
$ws.Range("T3").Value = 929
$ws.Range("U3").Value = -802
$ws.Range("V3").Value = 1066
$ws.Range("W3").Value = 0.96
$ws.Range("X3").Value = 0.74
$ws.Range("AA3").Value = 75.26000000000001
$ws.Range("AB3").Value = 1392.97
$ws.Range("AC3").Value = 1727
$ws.Range("AD3").Value = 100.19
$ws.Range("AE3").Value = 74153
$ws.Range("AF3").Value = 2.33
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 0.29
$ws.Range("AI3").Value = 28.96
$ws.Range("AJ3").Value = 3872480

# Row 4
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("D4").Value = 10690
$ws.Range("E4").Value = 214
$ws.Range("F4").Value = 214
$ws.Range("G4").Value = 181
$ws.Range("H4").Value = 143
$ws.Range("I4").Value = 143
$ws.Range("K4").Value = 5562
$ws.Range("L4").Value = 2590
$ws.Range("M4").Value = 2972
$ws.Range("N4").Value = 2972
$ws.Range("P4").Value = 194
$ws.Range("Q4").Value = 213
$ws.Range("R4").Value = -423
$ws.Range("S4").Value = 166
$ws.Range("T4").Value = 485
$ws.Range("U4").Value = -272
$ws.Range("V4").Value = 1251
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 1.34
$ws.Range("Y4").Value = 4.89
$ws.Range("Z4").Value = 2.7
$ws.Range("AA4").Value = 87.14
$ws.Range("AB4").Value = 1444.8
$ws.Range("AC4").Value = 3690
$ws.Range("AD4").Value = 39.97
$ws.Range("AE4").Value = 76745
$ws.Range("AF4").Value = 1.92
$ws.Range("AG4").Value = 600
$ws.Range("AH4").Value = 0.41
$ws.Range("AI4").Value = 16.26
$ws.Range("AJ4").Value = 3872480

# Row 5
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("D5").Value = 12075
$ws.Range("E5").Value = 298
$ws.Range("F5").Value = 298
$ws.Range("G5").Value = 253
$ws.Range("H5").Value = 206
$ws.Range("I5").Value = 206
$ws.Range("K5").Value = 5954
$ws.Range("L5").Value = 2803
$ws.Range("M5").Value = 3151
$ws.Range("N5").Value = 3151
$ws.Range("P5").Value = 194
$ws.Range("Q5").Value = 836
$ws.Range("R5").Value = -406
$ws.Range("S5").Value = 61
$ws.Range("T5").Value = 471
$ws.Range("U5").Value = 365
$ws.Range("V5").Value = 1335
$ws.Range("W5").Value = 2.47
$ws.Range("X5").Value = 1.71
$ws.Range("Y5").Value = 6.73
$ws.Range("Z5").Value = 3.58
$ws.Range("AA5").Value = 88.97
$ws.Range("AB5").Value = 1537.17
$ws.Range("AC5").Value = 5323
$ws.Range("AD5").Value = 24.89
$ws.Range("AE5").Value = 81364
$ws.Range("AF5").Value = 1.63
$ws.Range("AG5").Value = 750
$ws.Range("AH5").Value = 0.57
$ws.Range("AJ5").Value = 3872480

# Row 6
$ws.Range("D6").Value = 12786
$ws.Range("E6").Value = 274
$ws.Range("F6").Value = 274
$ws.Range("G6").Value = 135
$ws.Range("H6").Value = 85
$ws.Range("I6").Value = 85
$ws.Range("K6").Value = 6250
$ws.Range("L6").Value = 3081
$ws.Range("M6").Value = 3169
$ws.Range("N6").Value = 3169
$ws.Range("P6").Value = 194
$ws.Range("Q6").Value = 300
$ws.Range("R6").Value = -827
$ws.Range("S6").Value = 79
$ws.Range("T6").Value = 822
$ws.Range("U6").Value = -522
$ws.Range("V6").Value = 1444
$ws.Range("W6").Value = 2.14
$ws.Range("X6").Value = 0.67
$ws.Range("Y6").Value = 2.69
$ws.Range("Z6").Value = 1.4
$ws.Range("AA6").Value = 97.20999999999999
$ws.Range("AB6").Value = 1546.61
$ws.Range("AC6").Value = 2198
$ws.Range("AD6").Value = 38.85
$ws.Range("AE6").Value = 81839
$ws.Range("AF6").Value = 1.04
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 0.88
$ws.Range("AI6").Value = 34.12
$ws.Range("AJ6").Value = 3872480

# Row 7
$ws.Range("D7").Value = 13325
$ws.Range("E7").Value = 255
$ws.Range("G7").Value = 190
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 140
$ws.Range("K7").Value = 7060
$ws.Range("L7").Value = 3780
$ws.Range("M7").Value = 3280
$ws.Range("N7").Value = 3280
$ws.Range("P7").Value = 190
$ws.Range("Q7").Value = 760
$ws.Range("R7").Value = -750
$ws.Range("S7").Value = 30
$ws.Range("T7").Value = 720
$ws.Range("U7").Value = 40
$ws.Range("W7").Value = 1.91
$ws.Range("X7").Value = 1.13
$ws.Range("Y7").Value = 4.34
$ws.Range("Z7").Value = 2.25
$ws.Range("AA7").Value = 115.24
$ws.Range("AC7").Value = 3615
$ws.Range("AD7").Value = 18.64
$ws.Range("AE7").Value = 84704
$ws.Range("AF7").Value = 0.8
$ws.Range("AG7").Value = 750
$ws.Range("AH7").Value = 1.11
$ws.Range("AI7").Value = 20.75

# Row 8
$ws.Range("D8").Value = 14270
$ws.Range("E8").Value = 317
$ws.Range("G8").Value = 264
$ws.Range("H8").Value = 220
$ws.Range("I8").Value = 204
$ws.Range("K8").Value = 7320
$ws.Range("L8").Value = 3850
$ws.Range("M8").Value = 3470
$ws.Range("N8").Value = 3470
$ws.Range("P8").Value = 190
$ws.Range("Q8").Value = 650
$ws.Range("R8").Value = -580
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 550
$ws.Range("U8").Value = 100
$ws.Range("W8").Value = 2.22
$ws.Range("X8").Value = 1.54
$ws.Range("Y8").Value = 6.06
$ws.Range("Z8").Value = 3.06
$ws.Range("AA8").Value = 110.95
$ws.Range("AC8").Value = 5281
$ws.Range("AD8").Value = 12.76
$ws.Range("AE8").Value = 89611
$ws.Range("AF8").Value = 0.75
$ws.Range("AG8").Value = 750
$ws.Range("AH8").Value = 1.11
$ws.Range("AI8").Value = 14.2

# Row 9
$ws.Range("D9").Value = 15330
$ws.Range("E9").Value = 375
$ws.Range("G9").Value = 325
$ws.Range("H9").Value = 260
$ws.Range("I9").Value = 245
$ws.Range("K9").Value = 7620
$ws.Range("L9").Value = 3920
$ws.Range("M9").Value = 3700
$ws.Range("N9").Value = 3700
$ws.Range("P9").Value = 190
$ws.Range("Q9").Value = 700
$ws.Range("R9").Value = -590
$ws.Range("S9").Value = -10
$ws.Range("T9").Value = 550
$ws.Range("U9").Value = 150
$ws.Range("W9").Value = 2.45
$ws.Range("X9").Value = 1.7
$ws.Range("Y9").Value = 6.83
$ws.Range("Z9").Value = 3.48
$ws.Range("AA9").Value = 105.95
$ws.Range("AC9").Value = 6327
$ws.Range("AD9").Value = 10.65
$ws.Range("AE9").Value = 95551
$ws.Range("AF9").Value = 0.71
$ws.Range("AG9").Value = 750
$ws.Range("AH9").Value = 1.11
$ws.Range("AI9").Value = 11.86
